# Auto-generated Excel COM-interop script
# Refreshes the scraped cryptocurrency data on Sheet1:
#  - Price (D) and Volume(1h) (E) values are updated for many coins
#  - Hora (G) advances from "10" to "11" for every data row (2-51)
#  - Rows 15-24: "One" is newly listed and the following coins each
#    shift down one rank (Coin/Link/Price/Volume updated accordingly)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '''246.78'
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '''1.07%'
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(2, 7).Value = '''11'
$ws.Cells.Item(2, 7).Style = "Normal"
$ws.Cells.Item(3, 4).Value = '''30.17'
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = '''11.08%'
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Cells.Item(3, 7).Value = '''11'
$ws.Cells.Item(3, 7).Style = "Normal"
$ws.Cells.Item(4, 4).Value = '''5.172'
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '''0.29%'
$ws.Cells.Item(4, 5).Style = "Normal"
$ws.Cells.Item(4, 7).Value = '''11'
$ws.Cells.Item(4, 7).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '''1.81%'
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Cells.Item(5, 7).Value = '''11'
$ws.Cells.Item(5, 7).Style = "Normal"
$ws.Cells.Item(6, 4).Value = '''6.598'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '''1.86%'
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Cells.Item(6, 7).Value = '''11'
$ws.Cells.Item(6, 7).Style = "Normal"
$ws.Cells.Item(7, 4).Value = '''0.8583'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '''5.16%'
$ws.Cells.Item(7, 5).Style = "Normal"
$ws.Cells.Item(7, 7).Value = '''11'
$ws.Cells.Item(7, 7).Style = "Normal"
$ws.Cells.Item(8, 4).Value = '''0.8815'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '''6.10%'
$ws.Cells.Item(8, 5).Style = "Normal"
$ws.Cells.Item(8, 7).Value = '''11'
$ws.Cells.Item(8, 7).Style = "Normal"
$ws.Cells.Item(9, 4).Value = '''0.1368'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '''2.89%'
$ws.Cells.Item(9, 5).Style = "Normal"
$ws.Cells.Item(9, 7).Value = '''11'
$ws.Cells.Item(9, 7).Style = "Normal"
$ws.Cells.Item(10, 4).Value = '''0.07077'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '''2.36%'
$ws.Cells.Item(10, 5).Style = "Normal"
$ws.Cells.Item(10, 7).Value = '''11'
$ws.Cells.Item(10, 7).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '''-2.30%'
$ws.Cells.Item(11, 5).Style = "Normal"
$ws.Cells.Item(11, 7).Value = '''11'
$ws.Cells.Item(11, 7).Style = "Normal"
$ws.Cells.Item(12, 4).Value = '''0.09390'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '''0.07%'
$ws.Cells.Item(12, 5).Style = "Normal"
$ws.Cells.Item(12, 7).Value = '''11'
$ws.Cells.Item(12, 7).Style = "Normal"
$ws.Cells.Item(13, 4).Value = '''0.001510'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '''-0.76%'
$ws.Cells.Item(13, 5).Style = "Normal"
$ws.Cells.Item(13, 7).Value = '''11'
$ws.Cells.Item(13, 7).Style = "Normal"
$ws.Cells.Item(14, 4).Value = '''0.04142'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '''-2.14%'
$ws.Cells.Item(14, 5).Style = "Normal"
$ws.Cells.Item(14, 7).Value = '''11'
$ws.Cells.Item(14, 7).Style = "Normal"
$ws.Cells.Item(15, 2).Value = 'One'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Cells.Item(15, 4).Value = '''0.0005987'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '''-0.33%'
$ws.Cells.Item(15, 5).Style = "Normal"
$ws.Cells.Item(15, 7).Value = '''11'
$ws.Cells.Item(15, 7).Style = "Normal"
$ws.Cells.Item(16, 2).Value = 'TigerCash'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Cells.Item(16, 4).Value = '''0.006011'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '''-2.07%'
$ws.Cells.Item(16, 5).Style = "Normal"
$ws.Cells.Item(16, 7).Value = '''11'
$ws.Cells.Item(16, 7).Style = "Normal"
$ws.Cells.Item(17, 2).Value = 'LEO'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(17, 4).Value = '''3.495'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '''-2.01%'
$ws.Cells.Item(17, 5).Style = "Normal"
$ws.Cells.Item(17, 7).Value = '''11'
$ws.Cells.Item(17, 7).Style = "Normal"
$ws.Cells.Item(18, 2).Value = 'GateToken'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Cells.Item(18, 4).Value = '''3.069'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '''1.66%'
$ws.Cells.Item(18, 5).Style = "Normal"
$ws.Cells.Item(18, 7).Value = '''11'
$ws.Cells.Item(18, 7).Style = "Normal"
$ws.Cells.Item(19, 2).Value = 'BTSEToken'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Cells.Item(19, 4).Value = '''2.277'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '''-1.39%'
$ws.Cells.Item(19, 5).Style = "Normal"
$ws.Cells.Item(19, 7).Value = '''11'
$ws.Cells.Item(19, 7).Style = "Normal"
$ws.Cells.Item(20, 2).Value = 'BitpandaEcosystemToken'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Cells.Item(20, 4).Value = '''0.3185'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '''2.30%'
$ws.Cells.Item(20, 5).Style = "Normal"
$ws.Cells.Item(20, 7).Value = '''11'
$ws.Cells.Item(20, 7).Style = "Normal"
$ws.Cells.Item(21, 2).Value = 'LiechtensteinCryptoassetsExchange'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Cells.Item(21, 4).Value = '''0.03266'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '''5.72%'
$ws.Cells.Item(21, 5).Style = "Normal"
$ws.Cells.Item(21, 7).Value = '''11'
$ws.Cells.Item(21, 7).Style = "Normal"
$ws.Cells.Item(22, 2).Value = 'ProBitToken'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Cells.Item(22, 4).Value = '''0.1308'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '''1.26%'
$ws.Cells.Item(22, 5).Style = "Normal"
$ws.Cells.Item(22, 7).Value = '''11'
$ws.Cells.Item(22, 7).Style = "Normal"
$ws.Cells.Item(23, 2).Value = 'MCDex'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Cells.Item(23, 4).Value = '''3.531'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '''-5.56%'
$ws.Cells.Item(23, 5).Style = "Normal"
$ws.Cells.Item(23, 7).Value = '''11'
$ws.Cells.Item(23, 7).Style = "Normal"
$ws.Cells.Item(24, 2).Value = 'ZBToken'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Cells.Item(24, 4).Value = '''0.1379'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '''0.46%'
$ws.Cells.Item(24, 5).Style = "Normal"
$ws.Cells.Item(24, 7).Value = '''11'
$ws.Cells.Item(24, 7).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '''-1.00%'
$ws.Cells.Item(25, 5).Style = "Normal"
$ws.Cells.Item(25, 7).Value = '''11'
$ws.Cells.Item(25, 7).Style = "Normal"
$ws.Cells.Item(26, 4).Value = '''0.004501'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '''0.30%'
$ws.Cells.Item(26, 5).Style = "Normal"
$ws.Cells.Item(26, 7).Value = '''11'
$ws.Cells.Item(26, 7).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '''23.50%'
$ws.Cells.Item(27, 5).Style = "Normal"
$ws.Cells.Item(27, 7).Value = '''11'
$ws.Cells.Item(27, 7).Style = "Normal"
$ws.Cells.Item(28, 4).Value = '''0.0001383'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '''0.00%'
$ws.Cells.Item(28, 5).Style = "Normal"
$ws.Cells.Item(28, 7).Value = '''11'
$ws.Cells.Item(28, 7).Style = "Normal"
$ws.Cells.Item(29, 7).Value = '''11'
$ws.Cells.Item(29, 7).Style = "Normal"
$ws.Cells.Item(30, 7).Value = '''11'
$ws.Cells.Item(30, 7).Style = "Normal"
$ws.Cells.Item(31, 7).Value = '''11'
$ws.Cells.Item(31, 7).Style = "Normal"
$ws.Cells.Item(32, 7).Value = '''11'
$ws.Cells.Item(32, 7).Style = "Normal"
$ws.Cells.Item(33, 7).Value = '''11'
$ws.Cells.Item(33, 7).Style = "Normal"
$ws.Cells.Item(34, 7).Value = '''11'
$ws.Cells.Item(34, 7).Style = "Normal"
$ws.Cells.Item(35, 7).Value = '''11'
$ws.Cells.Item(35, 7).Style = "Normal"
$ws.Cells.Item(36, 7).Value = '''11'
$ws.Cells.Item(36, 7).Style = "Normal"
$ws.Cells.Item(37, 7).Value = '''11'
$ws.Cells.Item(37, 7).Style = "Normal"
$ws.Cells.Item(38, 7).Value = '''11'
$ws.Cells.Item(38, 7).Style = "Normal"
$ws.Cells.Item(39, 7).Value = '''11'
$ws.Cells.Item(39, 7).Style = "Normal"
$ws.Cells.Item(40, 4).Value = '''0.03793'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '''3.92%'
$ws.Cells.Item(40, 5).Style = "Normal"
$ws.Cells.Item(40, 7).Value = '''11'
$ws.Cells.Item(40, 7).Style = "Normal"
$ws.Cells.Item(41, 4).Value = '''0.005679'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '''-6.09%'
$ws.Cells.Item(41, 5).Style = "Normal"
$ws.Cells.Item(41, 7).Value = '''11'
$ws.Cells.Item(41, 7).Style = "Normal"
$ws.Cells.Item(42, 4).Value = '''0.1070'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '''1.58%'
$ws.Cells.Item(42, 5).Style = "Normal"
$ws.Cells.Item(42, 7).Value = '''11'
$ws.Cells.Item(42, 7).Style = "Normal"
$ws.Cells.Item(43, 4).Value = '''0.002199'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '''22.26%'
$ws.Cells.Item(43, 5).Style = "Normal"
$ws.Cells.Item(43, 7).Value = '''11'
$ws.Cells.Item(43, 7).Style = "Normal"
$ws.Cells.Item(44, 4).Value = '''0.01005'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '''24.09%'
$ws.Cells.Item(44, 5).Style = "Normal"
$ws.Cells.Item(44, 7).Value = '''11'
$ws.Cells.Item(44, 7).Style = "Normal"
$ws.Cells.Item(45, 4).Value = '''0.00005086'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '''-5.81%'
$ws.Cells.Item(45, 5).Style = "Normal"
$ws.Cells.Item(45, 7).Value = '''11'
$ws.Cells.Item(45, 7).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '''0.03%'
$ws.Cells.Item(46, 5).Style = "Normal"
$ws.Cells.Item(46, 7).Value = '''11'
$ws.Cells.Item(46, 7).Style = "Normal"
$ws.Cells.Item(47, 4).Value = '''0.08898'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '''-18.33%'
$ws.Cells.Item(47, 5).Style = "Normal"
$ws.Cells.Item(47, 7).Value = '''11'
$ws.Cells.Item(47, 7).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '''4.86%'
$ws.Cells.Item(48, 5).Style = "Normal"
$ws.Cells.Item(48, 7).Value = '''11'
$ws.Cells.Item(48, 7).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '''0.03%'
$ws.Cells.Item(49, 5).Style = "Normal"
$ws.Cells.Item(49, 7).Value = '''11'
$ws.Cells.Item(49, 7).Style = "Normal"
$ws.Cells.Item(50, 4).Value = '''0.0002000'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '''0.03%'
$ws.Cells.Item(50, 5).Style = "Normal"
$ws.Cells.Item(50, 7).Value = '''11'
$ws.Cells.Item(50, 7).Style = "Normal"
$ws.Cells.Item(51, 7).Value = '''11'
$ws.Cells.Item(51, 7).Style = "Normal"

Write-Host "Applied 139 cell updates to Sheet1"
